# Add PQ, ramp, limit trafo, variability function of time
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers for new columns D and E
$ws.Range("D1").Value = "CV_D"
$ws.Range("E1").Value = "CV_pv"

# Data values for D2:E25
$data = @(
    @(0.1, 0),
    @(0.1, 0),
    @(0.1, 0),
    @(0.2, 0),
    @(0.2, 0.1),
    @(0.2, 0.1),
    @(0.35, 0.1),
    @(0.35, 0.2),
    @(0.35, 0.2),
    @(0.35, 0.35),
    @(0.3, 0.35),
    @(0.3, 0.5),
    @(0.3, 0.35),
    @(0.3, 0.2),
    @(0.35, 0.2),
    @(0.35, 0.2),
    @(0.4, 0.1),
    @(0.45, 0.1),
    @(0.5, 0.1),
    @(0.55000000000000004, 0.1),
    @(0.5, 0),
    @(0.4, 0),
    @(0.3, 0),
    @(0.3, 0)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 4).Value = $data[$i][0]
    $ws.Cells.Item($row, 5).Value = $data[$i][1]
}

# Update the selection to match the recorded cursor position in the diff
$ws.Range("D23").Select()
